$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

# --- Sheet "About" ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: $newVersion"

$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Alardinskaya Coal Mine, Russia, M0767, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- Sheet "Boundaries and methane sources" ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 18; $row++) {
    $cell = $wsData.Range("S$row")
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value2 = $newVersion
    }
}
